# Read excel and list values
# This script reproduces the commit that:
#  - inserts two "section header" rows (K3 / K4) into the "sites" sheet,
#    pushing the existing data rows down and trimming now-unused trailing
#    '#' marker cells, and adds a trailing blank row
#  - adds a new worksheet ("Sheet1") that lists a couple of the values
#    (K4 header + the 046G/047G site codes)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "sites" sheet

# --- Insert the two new header rows -------------------------------------
# Row 1 becomes blank (will hold "K3"); what was row 1 slides to row 2.
$ws.Rows.Item(1).Insert()
# Row 3 becomes blank (will hold "K4"); what was row 2 slides to row 4.
$ws.Rows.Item(3).Insert()

# Fill K4 in first so it lands before K3 in the shared string table
# (matches the original authoring order of the workbook).
$ws.Range("A3").Value = "K4"
$ws.Range("A1").Value = "K3"

# Row insert cloned formatting into B:D of the two new rows - drop it so
# only column A carries a value/style on these header rows.
$ws.Range("B1:D1").Clear()
$ws.Range("B3:D3").Clear()

# The old trailing "#" marker cells are no longer part of the table.
$ws.Range("D2").Clear()
$ws.Range("D4").Clear()
$ws.Range("E5").Clear()

# Add the trailing blank (but styled) row.
$ws.Range("A6").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A6").VerticalAlignment = -4108     # xlCenter
$ws.Range("A6").Font.Bold = $false

# --- Add the new worksheet ----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Sheet1"

$newSheet.Range("A2").Value = "K4"
$newSheet.Range("A3").Value = "046G"
$newSheet.Range("B3").Value = "047G"

$newSheet.Range("A2").HorizontalAlignment = -4108
$newSheet.Range("A2").VerticalAlignment = -4108
$newSheet.Range("A3").HorizontalAlignment = -4108
$newSheet.Range("A3").VerticalAlignment = -4108
$newSheet.Range("B3").HorizontalAlignment = -4108
$newSheet.Range("B3").VerticalAlignment = -4108

# Leave the selection on the new sheet where the author left it, then
# re-activate "sites" so it stays the selected tab.
$newSheet.Range("D39").Select() | Out-Null
$ws.Range("B9").Select() | Out-Null
$ws.Activate() | Out-Null
